$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handback
#
# For both the "zh-cn" and "de-de" sheets, row 5 (the
# 42491a3d-c3ef-491f-b199-bdebb31e602c handback entry) gets its
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# / "Error Detail" columns filled in, and columns I, J, P are widened to 40
# (matching the other wide columns) to accommodate the new content.
# ---------------------------------------------------------------------------

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b98d2e51fdb15f74c38aef594adc8278bf21fd4f/e2e/42491a3d-c3ef-491f-b199-bdebb31e602c.md"
$targetDisplay = "42491a3d-c3ef-491f-b199-bdebb31e602c.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4f3c491c9070297ace39aa832bd94371fe16aa4e/e2e/42491a3d-c3ef-491f-b199-bdebb31e602c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b98d2e51fdb15f74c38aef594adc8278bf21fd4f/e2e/42491a3d-c3ef-491f-b199-bdebb31e602c.md."

$sheetInfo = @{
    "zh-cn" = @{
        HandbackFile = "42491a3d-c3ef-491f-b199-bdebb31e602c.a8249c8c0d458e720afe7d89222774018be778ae.zh-cn.xlf"
        HandbackDate = "2016-09-07 05:27:35"
    }
    "de-de" = @{
        HandbackFile = "42491a3d-c3ef-491f-b199-bdebb31e602c.a8249c8c0d458e720afe7d89222774018be778ae.de-de.xlf"
        HandbackDate = "2016-09-07 05:27:52"
    }
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $info = $sheetInfo[$sheetName]

    # Widen columns I (Latest Target File), J (Latest Handback File) and
    # P (Error Detail) to match the other wide (40-char) columns.
    $ws.Range("I1").ColumnWidth = $ws.Range("A1").ColumnWidth
    $ws.Range("J1").ColumnWidth = $ws.Range("A1").ColumnWidth
    $ws.Range("P1").ColumnWidth = $ws.Range("A1").ColumnWidth

    # Remember the pre-existing hyperlinks (A2:A6) so they can be re-added
    # in row order once the new I5 hyperlink is inserted between A5 and A6.
    $existingLinks = @()
    foreach ($hl in $ws.Hyperlinks) {
        $existingLinks += ,@($hl.Range.Address(), $hl.Address, $hl.TextToDisplay)
    }
    $ws.Hyperlinks.Delete()

    foreach ($link in $existingLinks) {
        $addr = ($link[0] -replace '\$', '')
        if ($addr -eq "A6") {
            # Insert the new "Latest Target File" hyperlink right before A6,
            # so the resulting order is A2, A3, A4, A5, I5, A6.
            $ws.Hyperlinks.Add($ws.Range("I5"), $targetUrl, $null, $null, $targetDisplay)
        }
        $ws.Hyperlinks.Add($ws.Range($addr), $link[1], $null, $null, $link[2])
    }

    # Row 5: fill in the handback results for 42491a3d-c3ef-491f-b199-bdebb31e602c
    $ws.Range("J5").Value = $info.HandbackFile
    $ws.Range("K5").Value = $info.HandbackDate
    $ws.Range("P5").Value = $errorDetail
}
